$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1100
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = $null
$ws.Range("H33").Value = 191.75
$ws.Range("I33").Value = 122.44444
$ws.Range("J33").Value = 399.66666
$ws.Range("K33").Value = 122.44444
$ws.Range("L33").Value = 399.66666
$ws.Range("M33").Value = 106.55556
$ws.Range("N33").Value = -857.66666
$ws.Range("H58").Value = 956.1667
$ws.Range("J58").Value = 1833
$ws.Range("L58").Value = 5499
$ws.Range("N58").Value = -5799
$ws.Range("H62").Value = 7938
$ws.Range("J62").Value = 7938
$ws.Range("L62").Value = 7938
$ws.Range("N62").Value = -9186
$ws.Range("H64").Value = 4579.8
$ws.Range("I64").Value = 3966.3333
$ws.Range("J64").Value = 5500
$ws.Range("K64").Value = 3966.3333
$ws.Range("L64").Value = 5500
$ws.Range("M64").Value = -3718.3333
$ws.Range("N64").Value = -5996
$ws.Range("H65").Value = 7938
$ws.Range("J65").Value = 7938
$ws.Range("L65").Value = 39690
$ws.Range("N65").Value = -45930
$ws.Range("H67").Value = 4579.8
$ws.Range("I67").Value = 3966.3333
$ws.Range("J67").Value = 5500
$ws.Range("K67").Value = 3966.3333
$ws.Range("L67").Value = 5500
$ws.Range("M67").Value = -3108.3333
$ws.Range("N67").Value = -7216
$ws.Range("H98").Value = 2620.2856
$ws.Range("I98").Value = 749.25
$ws.Range("K98").Value = 749.25
$ws.Range("M98").Value = 748.75
$ws.Range("H122").Value = 2620.2856
$ws.Range("I122").Value = 749.25
$ws.Range("K122").Value = 2247.75
$ws.Range("M122").Value = 202.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5505.706
$ws.Range("I32").Value = 5505.706
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5505.706
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -5218.706
$ws.Range("N32").Value = $null
$ws.Range("H61").Value = 3494.6
$ws.Range("I61").Value = 1735.9231
$ws.Range("K61").Value = 1735.9231
$ws.Range("M61").Value = -1523.9231
$ws.Range("H88").Value = 715.5
$ws.Range("I88").Value = 773
$ws.Range("J88").Value = 686.75
$ws.Range("K88").Value = 773
$ws.Range("L88").Value = 686.75
$ws.Range("M88").Value = -367
$ws.Range("N88").Value = -1498.75
$ws.Range("H91").Value = 715.5
$ws.Range("I91").Value = 773
$ws.Range("J91").Value = 686.75
$ws.Range("K91").Value = 773
$ws.Range("L91").Value = 686.75
$ws.Range("M91").Value = 631
$ws.Range("N91").Value = -3494.75
$ws.Range("H122").Value = 1552
$ws.Range("I122").Value = 1362.4
$ws.Range("K122").Value = 4087.2
$ws.Range("M122").Value = -1637.2
$ws.Range("H132").Value = 2897.8235
$ws.Range("I132").Value = 2556.2727
$ws.Range("K132").Value = 7668.8181
$ws.Range("M132").Value = -5138.8181
$ws.Range("H136").Value = 3494.6
$ws.Range("I136").Value = 1735.9231
$ws.Range("K136").Value = 5207.7693
$ws.Range("M136").Value = -2657.7693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 517.7692
$ws.Range("I22").Value = 448.77777
$ws.Range("J22").Value = 673
$ws.Range("K22").Value = 448.77777
$ws.Range("L22").Value = 673
$ws.Range("M22").Value = -275.77777
$ws.Range("N22").Value = -1019
$ws.Range("H86").Value = 4078.25
$ws.Range("I86").Value = 4268.6665
$ws.Range("J86").Value = 3507
$ws.Range("K86").Value = 4268.6665
$ws.Range("L86").Value = 3507
$ws.Range("M86").Value = -3145.6665
$ws.Range("N86").Value = -5753
$ws.Range("H89").Value = 4078.25
$ws.Range("I89").Value = 4268.6665
$ws.Range("J89").Value = 3507
$ws.Range("K89").Value = 21343.3325
$ws.Range("L89").Value = 17535
$ws.Range("M89").Value = -15727.3325
$ws.Range("N89").Value = -28767
$ws.Range("H105").Value = 3180
$ws.Range("I105").Value = 3457
$ws.Range("K105").Value = 3457
$ws.Range("M105").Value = -1710

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -650
$ws.Range("N22").Value = $null
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = $null
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").Value = $null
$ws.Range("H99").Value = 6193.6
$ws.Range("I99").Value = 4990
$ws.Range("J99").Value = 7999
$ws.Range("K99").Value = 4990
$ws.Range("L99").Value = 7999
$ws.Range("M99").Value = -3492
$ws.Range("N99").Value = -10995
$ws.Range("H126").Value = 6193.6
$ws.Range("I126").Value = 4990
$ws.Range("J126").Value = 7999
$ws.Range("K126").Value = 14970
$ws.Range("L126").Value = 23997
$ws.Range("M126").Value = -12500
$ws.Range("N126").Value = -28937
$ws.Range("H134").Value = 2281.2144
$ws.Range("I134").Value = 2476.9167
$ws.Range("K134").Value = 7430.750100000001
$ws.Range("M134").Value = -4895.750100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 79966.336
$ws.Range("J37").Value = 79966.336
$ws.Range("L37").Value = 239899.008
$ws.Range("N37").Value = -240123.008
$ws.Range("H64").Value = 17945
$ws.Range("J64").Value = 17945
$ws.Range("L64").Value = 53835
$ws.Range("N64").Value = -54375
$ws.Range("H67").Value = 17945
$ws.Range("J67").Value = 17945
$ws.Range("L67").Value = 53835
$ws.Range("N67").Value = -55707
$ws.Range("H95").Value = 5999
$ws.Range("J95").Value = 5999
$ws.Range("L95").Value = 17997
$ws.Range("N95").Value = -22115

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 173.5
$ws.Range("I2").Value = 209
$ws.Range("J2").Value = 19.666666
$ws.Range("K2").Value = 209
$ws.Range("L2").Value = 19.666666
$ws.Range("M2").Value = -96
$ws.Range("N2").Value = -245.666666
$ws.Range("H34").Value = 46500
$ws.Range("J34").Value = 46500
$ws.Range("L34").Value = 46500
$ws.Range("N34").Value = -47036
$ws.Range("H76").Value = 46500
$ws.Range("J76").Value = 46500
$ws.Range("L76").Value = 46500
$ws.Range("N76").Value = -47130
$ws.Range("H79").Value = 46500
$ws.Range("J79").Value = 46500
$ws.Range("L79").Value = 46500
$ws.Range("N79").Value = -48684
$ws.Range("H104").Value = 65994.5
$ws.Range("J104").Value = 65994.5
$ws.Range("L104").Value = 65994.5
$ws.Range("N104").Value = -72982.5
$ws.Range("H122").Value = 2194.6667
$ws.Range("I122").Value = 2194.6667
$ws.Range("K122").Value = 6584.000100000001
$ws.Range("M122").Value = -4134.000100000001
$ws.Range("H132").Value = 1497.8
$ws.Range("I132").Value = 1497
$ws.Range("K132").Value = 4491
$ws.Range("M132").Value = -1961

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2951
$ws.Range("I22").Value = 2898.2856
$ws.Range("J22").Value = 3043.25
$ws.Range("K22").Value = 2898.2856
$ws.Range("L22").Value = 3043.25
$ws.Range("M22").Value = -2603.2856
$ws.Range("N22").Value = -3633.25
$ws.Range("H27").Value = 2951
$ws.Range("I27").Value = 2898.2856
$ws.Range("J27").Value = 3043.25
$ws.Range("K27").Value = 2898.2856
$ws.Range("L27").Value = 3043.25
$ws.Range("M27").Value = -2791.2856
$ws.Range("N27").Value = -3257.25
$ws.Range("H40").Value = 3370.7273
$ws.Range("I40").Value = 1995.8
$ws.Range("J40").Value = 4516.5
$ws.Range("K40").Value = 1995.8
$ws.Range("L40").Value = 4516.5
$ws.Range("M40").Value = -1859.8
$ws.Range("N40").Value = -4788.5
$ws.Range("H46").Value = 65170.125
$ws.Range("I46").Value = 74051.57000000001
$ws.Range("K46").Value = 74051.57000000001
$ws.Range("M46").Value = -73863.57000000001
$ws.Range("H74").Value = 69729
$ws.Range("I74").Value = 69729
$ws.Range("K74").Value = 69729
$ws.Range("M74").Value = -68731
$ws.Range("H77").Value = 69729
$ws.Range("I77").Value = 69729
$ws.Range("K77").Value = 209187
$ws.Range("M77").Value = -204195
$ws.Range("H82").Value = 1109.4286
$ws.Range("I82").Value = 1109.4286
$ws.Range("K82").Value = 1109.4286
$ws.Range("M82").Value = -748.4286
$ws.Range("H85").Value = 1109.4286
$ws.Range("I85").Value = 1109.4286
$ws.Range("K85").Value = 1109.4286
$ws.Range("M85").Value = 138.5714
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = $null
$ws.Range("H122").Value = 7406.654
$ws.Range("I122").Value = 7965.2666
$ws.Range("K122").Value = 23895.7998
$ws.Range("M122").Value = -21445.7998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").Value = $null
$ws.Range("H81").Value = 3334216.8
$ws.Range("I81").Value = 1324
$ws.Range("K81").Value = 2648
$ws.Range("M81").Value = -1587
$ws.Range("H84").Value = 3334216.8
$ws.Range("I84").Value = 1324
$ws.Range("K84").Value = 13240
$ws.Range("M84").Value = -7936
$ws.Range("H100").Value = 6252256
$ws.Range("I100").Value = 7144506.5
$ws.Range("K100").Value = 14289013
$ws.Range("M100").Value = -14288472
$ws.Range("H122").Value = 1054.909
$ws.Range("I122").Value = 901
$ws.Range("K122").Value = 2703
$ws.Range("M122").Value = -253
